# RDCC-5182 Added Version check
# Add a new "VERSION" worksheet after the existing "Service to CW Roles
# Mapping" sheet, and populate it with a "File version" / "vx.xx" pair
# at row 6 (matching the committed workbook layout).

$wb = $excel.ActiveWorkbook

$mappingSheet = $wb.Worksheets.Item(1)

# Insert the new sheet right after the mapping sheet.
$versionSheet = $wb.Worksheets.Add($null, $mappingSheet)
$versionSheet.Name = "VERSION"

$versionSheet.Range("A6").Value = "File version"
$versionSheet.Range("B6").Value = "vx.xx"

# Make the new sheet the active tab/selection, like the authored workbook.
[void]$versionSheet.Range("B6").Select()
